$wb = $excel.ActiveWorkbook

# Rename the five "Include from Terminologia de ..." sheets to "Include #N"
$names = @("Include #0", "Include #1", "Include #2", "Include #3", "Include #4")
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 2)
    $ws.Name = $names[$i]
}

# Update the Contact value on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B10").Value = "null (http://www.saude.gov.br)"
